$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 67; this shifts existing rows 67-72 down to 68-73
$ws.Rows.Item(67).Insert()

# Populate the new row 67 with the new data record
$ws.Cells.Item(67, 1).Value = 4
$ws.Cells.Item(67, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(67, 3).Value = "Los Lagos"
$ws.Cells.Item(67, 4).Value = 45244
$ws.Cells.Item(67, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(67, 5).Value = 10
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100101
$ws.Cells.Item(67, 8).Value = "Berries"
$ws.Cells.Item(67, 9).Value = 100101001
$ws.Cells.Item(67, 10).Value = "Arándano (blue)"
$ws.Cells.Item(67, 11).Value = "Sin especificar"
$ws.Cells.Item(67, 12).Value = "Primera"
$ws.Cells.Item(67, 13).Value = 50
$ws.Cells.Item(67, 14).Value = 10500
$ws.Cells.Item(67, 15).Value = 10500
$ws.Cells.Item(67, 16).Value = 10500
$ws.Cells.Item(67, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(67, 18).Value = "Región del Maule"
$ws.Cells.Item(67, 19).Value = 5250
$ws.Cells.Item(67, 20).Value = 2
